# Auto-generated: apply scheduled market-data refresh values to Durandal_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 53
$ws.Range("H53").Value = 846.2778
$ws.Range("I53").Value = 913.9231
$ws.Range("K53").Value = 913.9231
$ws.Range("M53").Value = -276.9231
# row 64
$ws.Range("H64").Value = 3782.2942
$ws.Range("I64").Value = 3667.75
$ws.Range("J64").Value = 3884.111
$ws.Range("K64").Value = 3667.75
$ws.Range("L64").Value = 3884.111
$ws.Range("M64").Value = -3419.75
$ws.Range("N64").Value = -4380.111
# row 67
$ws.Range("H67").Value = 3782.2942
$ws.Range("I67").Value = 3667.75
$ws.Range("J67").Value = 3884.111
$ws.Range("K67").Value = 3667.75
$ws.Range("L67").Value = 3884.111
$ws.Range("M67").Value = -2809.75
$ws.Range("N67").Value = -5600.111
# row 74
$ws.Range("H74").Value = 4477.3335
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 4966
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 4966
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -6838
# row 76
$ws.Range("H76").Value = 2649054.5
$ws.Range("I76").Value = 2852476.5
$ws.Range("K76").Value = 2852476.5
$ws.Range("M76").Value = -2852161.5
# row 77
$ws.Range("H77").Value = 4477.3335
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 4966
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 24830
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -34190
# row 79
$ws.Range("H79").Value = 2649054.5
$ws.Range("I79").Value = 2852476.5
$ws.Range("K79").Value = 2852476.5
$ws.Range("M79").Value = -2851384.5
# row 116
$ws.Range("H116").Value = 6668.385
$ws.Range("I116").Value = 8860.3125
$ws.Range("J116").Value = 3161.3
$ws.Range("K116").Value = 8860.3125
$ws.Range("L116").Value = 3161.3
$ws.Range("M116").Value = -5418.3125
$ws.Range("N116").Value = -10045.3

$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 2274.2727
$ws.Range("I45").Value = 2182.6128
$ws.Range("J45").Value = 2492.8462
$ws.Range("K45").Value = 2182.6128
$ws.Range("L45").Value = 2492.8462
$ws.Range("M45").Value = -1805.6128
$ws.Range("N45").Value = -3246.8462
# row 61
$ws.Range("H61").Value = 2792.85
$ws.Range("I61").Value = 2791.2778
$ws.Range("J61").Value = 2807
$ws.Range("K61").Value = 2791.2778
$ws.Range("L61").Value = 2807
$ws.Range("M61").Value = -2579.2778
$ws.Range("N61").Value = -3231
# row 63
$ws.Range("H63").Value = 4146.875
$ws.Range("I63").Value = 1481.25
$ws.Range("J63").Value = 6812.5
$ws.Range("K63").Value = 1481.25
$ws.Range("L63").Value = 6812.5
$ws.Range("M63").Value = -795.25
$ws.Range("N63").Value = -8184.5
# row 66
$ws.Range("H66").Value = 4146.875
$ws.Range("I66").Value = 1481.25
$ws.Range("J66").Value = 6812.5
$ws.Range("K66").Value = 7406.25
$ws.Range("L66").Value = 34062.5
$ws.Range("M66").Value = -3974.25
$ws.Range("N66").Value = -40926.5
# row 136
$ws.Range("H136").Value = 2792.85
$ws.Range("I136").Value = 2791.2778
$ws.Range("J136").Value = 2807
$ws.Range("K136").Value = 8373.8334
$ws.Range("L136").Value = 8421
$ws.Range("M136").Value = -5823.8334
$ws.Range("N136").Value = -13521

$ws = $wb.Worksheets.Item("BSM")
# row 19
$ws.Range("H19").Value = 9000
$ws.Range("J19").Value = 9000
$ws.Range("L19").Value = 9000
$ws.Range("N19").Value = -9346
# row 20
$ws.Range("H20").Value = 3110.0588
$ws.Range("I20").Value = 2559.7646
$ws.Range("J20").Value = 3660.353
$ws.Range("K20").Value = 2559.7646
$ws.Range("L20").Value = 3660.353
$ws.Range("M20").Value = -2312.7646
$ws.Range("N20").Value = -4154.353
# row 105
$ws.Range("H105").Value = 1912
$ws.Range("I105").Value = 1918.7188
$ws.Range("J105").Value = 1840.3334
$ws.Range("K105").Value = 1918.7188
$ws.Range("L105").Value = 1840.3334
$ws.Range("M105").Value = -171.7188000000001
$ws.Range("N105").Value = -5334.3334
# row 134
$ws.Range("H134").Value = 7047.174
$ws.Range("I134").Value = 740.6923
$ws.Range("J134").Value = 15245.6
$ws.Range("K134").Value = 2222.0769
$ws.Range("L134").Value = 45736.8
$ws.Range("M134").Value = 312.9231
$ws.Range("N134").Value = -50806.8

$ws = $wb.Worksheets.Item("CRP")
# row 25
$ws.Range("H25").Value = 111
$ws.Range("I25").Value = 111
$ws.Range("K25").Value = 111
$ws.Range("M25").Value = 63
# row 31
$ws.Range("H31").Value = 12186.6455
$ws.Range("I31").Value = 4926.7827
$ws.Range("J31").Value = 16468.104
$ws.Range("K31").Value = 4926.7827
$ws.Range("L31").Value = 16468.104
$ws.Range("M31").Value = -4631.7827
$ws.Range("N31").Value = -17058.104
# row 34
$ws.Range("H34").Value = 12186.6455
$ws.Range("I34").Value = 4926.7827
$ws.Range("J34").Value = 16468.104
$ws.Range("K34").Value = 4926.7827
$ws.Range("L34").Value = 16468.104
$ws.Range("M34").Value = -4724.7827
$ws.Range("N34").Value = -16872.104
# row 62
$ws.Range("H62").Value = 3830.875
$ws.Range("I62").Value = 3536.25
$ws.Range("J62").Value = 4125.5
$ws.Range("K62").Value = 3536.25
$ws.Range("L62").Value = 4125.5
$ws.Range("M62").Value = -2912.25
$ws.Range("N62").Value = -5373.5
# row 65
$ws.Range("H65").Value = 3830.875
$ws.Range("I65").Value = 3536.25
$ws.Range("J65").Value = 4125.5
$ws.Range("K65").Value = 17681.25
$ws.Range("L65").Value = 20627.5
$ws.Range("M65").Value = -14561.25
$ws.Range("N65").Value = -26867.5
# row 94
$ws.Range("H94").Value = 5070
$ws.Range("I94").Value = 1467.5
$ws.Range("J94").Value = 5970.625
$ws.Range("K94").Value = 1467.5
$ws.Range("L94").Value = 5970.625
$ws.Range("M94").Value = -1016.5
$ws.Range("N94").Value = -6872.625
# row 134
$ws.Range("H134").Value = 545.4595
$ws.Range("I134").Value = 527.2083
$ws.Range("J134").Value = 579.1539
$ws.Range("K134").Value = 1581.6249
$ws.Range("L134").Value = 1737.4617
$ws.Range("M134").Value = 953.3751
$ws.Range("N134").Value = -6807.4617

$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 20834724
$ws.Range("I131").Value = 3633.3333
$ws.Range("J131").Value = 23810594
$ws.Range("K131").Value = 10899.9999
$ws.Range("L131").Value = 71431782
$ws.Range("M131").Value = -5859.999899999999
$ws.Range("N131").Value = -71441862

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 14810154
$ws.Range("I70").Value = 21641488
$ws.Range("J70").Value = 8930.75
$ws.Range("K70").Value = 21641488
$ws.Range("L70").Value = 8930.75
$ws.Range("M70").Value = -21641218
$ws.Range("N70").Value = -9470.75
# row 73
$ws.Range("H73").Value = 14810154
$ws.Range("I73").Value = 21641488
$ws.Range("J73").Value = 8930.75
$ws.Range("K73").Value = 21641488
$ws.Range("L73").Value = 8930.75
$ws.Range("M73").Value = -21640552
$ws.Range("N73").Value = -10802.75
# row 80
$ws.Range("H80").Value = 103217.55
$ws.Range("I80").Value = 113139.3
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 113139.3
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -112141.3
$ws.Range("N80").Value = -5996
# row 83
$ws.Range("H83").Value = 103217.55
$ws.Range("I83").Value = 113139.3
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 565696.5
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -560704.5
$ws.Range("N83").Value = -29984
# row 113
$ws.Range("H113").Value = 6472347.5
$ws.Range("I113").Value = 20001474
$ws.Range("J113").Value = 835211.8
$ws.Range("K113").Value = 20001474
$ws.Range("L113").Value = 835211.8
$ws.Range("M113").Value = -19999304
$ws.Range("N113").Value = -839551.8

$ws = $wb.Worksheets.Item("LTW")
# row 136
$ws.Range("H136").Value = 3664.4866
$ws.Range("I136").Value = 2389.0952
$ws.Range("K136").Value = 7167.285600000001
$ws.Range("M136").Value = -4617.285600000001

$ws = $wb.Worksheets.Item("WVR")
# row 2
$ws.Range("H2").Value = 25004990
$ws.Range("I2").Value = 25004990
$ws.Range("K2").Value = 25004990
$ws.Range("M2").Value = -25004878
# row 4
$ws.Range("H4").Value = 20001396
$ws.Range("I4").Value = 33333728
$ws.Range("J4").Value = 2901.5
$ws.Range("K4").Value = 33333728
$ws.Range("L4").Value = 2901.5
$ws.Range("M4").Value = -33333615
$ws.Range("N4").Value = -3127.5
# row 6
$ws.Range("H6").Value = 10423.5
$ws.Range("I6").Value = 341.875
$ws.Range("J6").Value = 50750
$ws.Range("K6").Value = 341.875
$ws.Range("L6").Value = 50750
$ws.Range("M6").Value = -226.875
$ws.Range("N6").Value = -50980
# row 15
$ws.Range("H15").Value = 6801.5
$ws.Range("I15").Value = 6206
$ws.Range("K15").Value = 6206
$ws.Range("M15").Value = -5918
# row 29
$ws.Range("H29").Value = 650
$ws.Range("I29").Value = 650
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 650
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -360
$ws.Range("N29").ClearContents() | Out-Null
